$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 0. Build a scratch single-run "Git, IntelliJ" FormattedText blob at
#    the very end of the document. We will use it later (step 4) to
#    replace the existing two-run "Git" / ", IntelliJ" pair with a
#    single merged run, without Word's adjacent-run auto-merge logic
#    bleeding into the following ", Android Studio, Eclipse" run.
# -----------------------------------------------------------------
$content = $d.Content
$endRange = $d.Range($content.End - 1, $content.End - 1)
$endRange.InsertParagraphAfter()
$scratchIdx = $d.Paragraphs.Count
$scratchRange = $d.Paragraphs($scratchIdx).Range
$scratchRange.InsertAfter("PLACEHOLDER")
$scratchRange2 = $d.Paragraphs($scratchIdx).Range
$scratchRange2.Find.Execute("PLACEHOLDER", $false, $false, $false, $false, $false, $true, 1, $false, "Git, IntelliJ", 2)
$scratchRange3 = $d.Paragraphs($scratchIdx).Range
$gitIntelliJFT = $scratchRange3.FormattedText

# -----------------------------------------------------------------
# 1. The "Asteroids Interactive Single Player Game using Python"
#    project bullet moves down one bullet (becomes the "Population
#    growth..." bullet), and the first bullet becomes the new
#    "Email Spam classifier..." text.
# -----------------------------------------------------------------
$r26 = $d.Paragraphs(26).Range
$projFT = $r26.FormattedText

$r27 = $d.Paragraphs(27).Range
$r27.FormattedText = $projFT

$r26b = $d.Paragraphs(26).Range
$r26b.Find.Execute("Asteroids Interactive Single Player Game using Python", $false, $false, $false, $false, $false, $true, 1, $false, "Email Spam classifier and Hand Written Digit Recognition using Matlab (Supervised Learning, ML)", 2)

# -----------------------------------------------------------------
# 2. Move the "_GoBack" bookmark from the Git/IntelliJ bullet to the
#    end of the "Ruzzle Game Solver using TRIE datastructure" bullet.
# -----------------------------------------------------------------
# 2a. Remove the bookmark from its old location first (so there is no
#     name clash when we re-add it elsewhere).
$gb = $d.Bookmarks("_GoBack")
$gb.Delete()

# 2b. Insert a throwaway marker character at the end of the
#     "...datastructure" paragraph, wrap it with the bookmark, then
#     delete the marker text again - this leaves a clean, empty
#     bookmark sitting right after "structure" (adding a bookmark to
#     a truly empty/collapsed range at that exact boundary is not
#     reliable, so we use a marker character as a safe stand-in).
$p28 = $d.Paragraphs(28).Range
$markerPos = $d.Range($p28.End - 1, $p28.End - 1)
$markerPos.InsertAfter("X")
$p28b = $d.Paragraphs(28).Range
$markerRange = $d.Range($p28b.End - 2, $p28b.End - 1)
$d.Bookmarks.Add("_GoBack", $markerRange)
$markerRange2 = $d.Range($p28b.End - 2, $p28b.End - 1)
$markerRange2.Text = ""

# -----------------------------------------------------------------
# 3. Merge the "Git" and ", IntelliJ" runs into a single run reading
#    "Git, IntelliJ" (leaving the following ", Android Studio,
#    Eclipse" run untouched).
# -----------------------------------------------------------------
$r35 = $d.Paragraphs(35).Range
$gitStart = $r35.Start
$gitSub = $d.Range($gitStart, $gitStart + 13)
$gitSub.FormattedText = $gitIntelliJFT

# -----------------------------------------------------------------
# 4. Clean up the scratch paragraph used to build the FormattedText.
# -----------------------------------------------------------------
$scratchFinal = $d.Paragraphs($d.Paragraphs.Count).Range
$scratchFinal.Delete()

Write-Output "done"
